$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.22080000000002
$ws.Range("A10").Value = -20.48809999999996
$ws.Range("A12").Value = -22.41750000000003
$ws.Range("A18").Value = -22.25650000000002
$ws.Range("A25").Value = -22.28480000000003
